$wb = $excel.ActiveWorkbook
$wb.Worksheets.Item("Group1").Name = "G1"
$wb.Worksheets.Item("Group2").Name = "G2"
$wb.Worksheets.Item("Sheet1").Activate()
